$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: financial period headers (D8:H8)
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# Row 9: publish dates (D9:H9)
$ws.Range("D9").Value = "1399-02-11 (8)"
$ws.Range("E9").Value = "1400-02-08 (8)"
$ws.Range("F9").Value = "1401-02-06 (9)"
$ws.Range("G9").Value = "1402-02-10 (8)"
$ws.Range("H9").Value = "1402-02-10 (2)"

# Data rows D:H (shift left one year + new rightmost year)
# Row 12
$ws.Range("D12").Value = 74216
$ws.Range("E12").Value = 197611
$ws.Range("F12").Value = 897732
$ws.Range("G12").Value = 838674
$ws.Range("H12").Value = 2486700

# Row 13
$ws.Range("D13").Value = -1763
$ws.Range("E13").Value = -37564
$ws.Range("F13").Value = -87438
$ws.Range("G13").Value = -118896
$ws.Range("H13").Value = -262122

# Row 14
$ws.Range("D14").Value = 72453
$ws.Range("E14").Value = 160047
$ws.Range("F14").Value = 810294
$ws.Range("G14").Value = 719778
$ws.Range("H14").Value = 2224578

# Row 16
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 1338
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0

# Row 17
$ws.Range("D17").Value = -1719
$ws.Range("E17").Value = -82120
$ws.Range("F17").Value = -123090
$ws.Range("G17").Value = -680103
$ws.Range("H17").Value = -1150359

# Row 18
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 6599
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0

# Row 19
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0

# Row 20
$ws.Range("D20").Value = -897
$ws.Range("E20").Value = -392
$ws.Range("F20").Value = -513
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0

# Row 21
$ws.Range("D21").Value = -39694
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0

# Row 22
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 34519
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0

# Row 23
$ws.Range("D23").Value = 83124
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# Row 24
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0

# Row 25
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 1765
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 49868

# Row 26
$ws.Range("D26").Value = -2
$ws.Range("E26").Value = -16
$ws.Range("F26").Value = -20
$ws.Range("G26").Value = -49934
$ws.Range("H26").Value = -78

# Row 27
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0

# Row 28
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0

# Row 29
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0

# Row 30
$ws.Range("D30").Value = 3360
$ws.Range("E30").Value = 166
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0

# Row 31
$ws.Range("D31").Value = 1260
$ws.Range("E31").Value = 262
$ws.Range("F31").Value = 888
$ws.Range("G31").Value = 9397
$ws.Range("H31").Value = 4678

# Row 32
$ws.Range("D32").Value = 45432
$ws.Range("E32").Value = -44478
$ws.Range("F32").Value = -116134
$ws.Range("G32").Value = -720640
$ws.Range("H32").Value = -1095891

# Row 33
$ws.Range("D33").Value = 117885
$ws.Range("E33").Value = 115569
$ws.Range("F33").Value = 694160
$ws.Range("G33").Value = -862
$ws.Range("H33").Value = 1128687

# Row 35
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 437490
$ws.Range("H35").Value = 37988

# Row 36
$ws.Range("D36").Value = "-"
$ws.Range("E36").Value = "-"
$ws.Range("F36").Value = "-"
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0

# Row 37
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0

# Row 38
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = -22000
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0

# Row 39
$ws.Range("D39").Value = 133850
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 600000
$ws.Range("H39").Value = 728989

# Row 40
$ws.Range("D40").Value = -114266
$ws.Range("E40").Value = -45573
$ws.Range("F40").Value = -152879
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = -638314

# Row 41
$ws.Range("D41").Value = -78607
$ws.Range("E41").Value = -16531
$ws.Range("F41").Value = -9111
$ws.Range("G41").Value = -252
$ws.Range("H41").Value = -59018

# Row 42
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0

# Row 43
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0

# Row 44
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0

# Row 45
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 0

# Row 46
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0

# Row 47
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0

# Row 48
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0

# Row 49
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0

# Row 50
$ws.Range("D50").Value = -40378
$ws.Range("E50").Value = -16829
$ws.Range("F50").Value = -310482
$ws.Range("G50").Value = -578781
$ws.Range("H50").Value = -1296297

# Row 51
$ws.Range("D51").Value = -99401
$ws.Range("E51").Value = -78933
$ws.Range("F51").Value = -494472
$ws.Range("G51").Value = 458457
$ws.Range("H51").Value = -1226652

# Row 52
$ws.Range("D52").Value = 18484
$ws.Range("E52").Value = 36636
$ws.Range("F52").Value = 199688
$ws.Range("G52").Value = 457595
$ws.Range("H52").Value = -97965

# Row 53
$ws.Range("D53").Value = 6611
$ws.Range("E53").Value = 25865
$ws.Range("F53").Value = 62627
$ws.Range("G53").Value = 262429
$ws.Range("H53").Value = 713802

# Row 54
$ws.Range("D54").Value = 770
$ws.Range("E54").Value = 126
$ws.Range("F54").Value = 114
$ws.Range("G54").Value = -6222
$ws.Range("H54").Value = 78

# Row 55
$ws.Range("D55").Value = 25865
$ws.Range("E55").Value = 62627
$ws.Range("F55").Value = 262429
$ws.Range("G55").Value = 713802
$ws.Range("H55").Value = 615915

# Row 56
$ws.Range("D56").Value = 5743
$ws.Range("E56").Value = 148300
$ws.Range("F56").Value = 140093
$ws.Range("G56").Value = 139118
$ws.Range("H56").Value = 0
